$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting existing data (D:K) to (F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting from the (now shifted) adjacent columns F and G
# into the newly inserted D and E columns so styles match exactly.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("G5:G102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D and E columns with the new quarter data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1724000
$ws.Range("E8").Value = 1818000
$ws.Range("D9").Value = 1311000
$ws.Range("E9").Value = 1367000
$ws.Range("D10").Value = 413000
$ws.Range("E10").Value = 451000
$ws.Range("D12").Value = 23000
$ws.Range("E12").Value = 21000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3000
$ws.Range("E14").Value = 8000
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 1502000
$ws.Range("E17").Value = 1563000
$ws.Range("D18").Value = 222000
$ws.Range("E18").Value = 255000
$ws.Range("D20").Value = -22000
$ws.Range("E20").Value = -27000
$ws.Range("D21").Value = 310000
$ws.Range("E21").Value = 335000
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 200000
$ws.Range("E23").Value = 228000
$ws.Range("D24").Value = 29000
$ws.Range("E24").Value = 67000
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 171000
$ws.Range("E26").Value = 161000
$ws.Range("D27").Value = 171000
$ws.Range("E27").Value = 161000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 22000
$ws.Range("E32").Value = 27000
$ws.Range("D33").Value = 171000
$ws.Range("E33").Value = 161000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 171000
$ws.Range("E35").Value = 161000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 78000
$ws.Range("E41").Value = 136000
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 794000
$ws.Range("E43").Value = 955000
$ws.Range("D44").Value = 1072000
$ws.Range("E44").Value = 1032000
$ws.Range("D45").Value = 76000
$ws.Range("E45").Value = 118000
$ws.Range("D46").Value = 2020000
$ws.Range("E46").Value = 2241000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 3811000
$ws.Range("E48").Value = 3741000
$ws.Range("D49").Value = 3728000
$ws.Range("E49").Value = 3759000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 212000
$ws.Range("E52").Value = 305000
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 9771000
$ws.Range("E54").Value = 10046000
$ws.Range("D57").Value = 842000
$ws.Range("E57").Value = 1379000
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = 5000
$ws.Range("D59").Value = 436000
$ws.Range("E59").Value = 2000
$ws.Range("D60").Value = 1278000
$ws.Range("E60").Value = 1386000
$ws.Range("D61").Value = 3362000
$ws.Range("E61").Value = 3669000
$ws.Range("D62").Value = 807000
$ws.Range("E62").Value = 785000
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 5488000
$ws.Range("E66").Value = 5881000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 2013000
$ws.Range("E72").Value = 1867000
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 4283000
$ws.Range("E76").Value = 4165000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 171000
$ws.Range("E81").Value = 161000
$ws.Range("D83").Value = 110000
$ws.Range("E83").Value = 107000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 297000
$ws.Range("E89").Value = 200000
$ws.Range("D91").Value = -112000
$ws.Range("E91").Value = -121000
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -38000
$ws.Range("E94").Value = -121000
$ws.Range("D96").Value = -22000
$ws.Range("E96").Value = -24000
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -306000
$ws.Range("E100").Value = -98000
$ws.Range("D101").Value = -12000
$ws.Range("E101").Value = 7000
$ws.Range("D102").Value = -59000
$ws.Range("E102").Value = -12000
